$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.833.15'
$ws.Range('E2').Value = '  +1.41%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.838.54'
$ws.Range('E3').Value = '  +1.51%  '

$ws.Range('E4').Value = '  +0.59%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.82'
$ws.Range('E5').Value = '  +1.08%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.42%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4657'
$ws.Range('E7').Value = '  +3.38%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3639'
$ws.Range('E8').Value = '  +1.58%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07112'
$ws.Range('E9').Value = '  +0.78%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9128'
$ws.Range('E10').Value = '  +2.63%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07693'
$ws.Range('E11').Value = '  -1.00%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.48'
$ws.Range('E12').Value = '  +0.78%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.860.82'
$ws.Range('E13').Value = '  +2.85%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.268'
$ws.Range('E14').Value = '  +0.04%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.380'
$ws.Range('E15').Value = '  +1.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.58'
$ws.Range('E16').Value = '  +3.10%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.51%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008570'
$ws.Range('E18').Value = '  +0.59%  '

$ws.Range('E19').Value = '  +0.43%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.866.11'
$ws.Range('E20').Value = '  +1.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.35'
$ws.Range('E21').Value = '  +1.30%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.011'
$ws.Range('E22').Value = '  +1.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.62'
$ws.Range('E23').Value = '  +1.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.931'
$ws.Range('E24').Value = '  -1.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.23'
$ws.Range('E25').Value = '  +0.81%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.18'
$ws.Range('E26').Value = '  +2.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.017'
$ws.Range('E27').Value = '  -1.67%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '113.87'
$ws.Range('E28').Value = '  +1.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.873'
$ws.Range('E29').Value = '  +0.80%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08856'
$ws.Range('E30').Value = '  +1.99%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.202'
$ws.Range('E31').Value = '  +2.12%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.167'
$ws.Range('E32').Value = '  +5.45%  '

$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.789'
$ws.Range('E33').Value = '  +1.57%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7411'
$ws.Range('E34').Value = '  -0.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.450'
$ws.Range('E35').Value = '  +0.43%  '

$ws.Range('E36').Value = '  +1.02%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01931'
$ws.Range('E37').Value = '  +0.25%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.961'
$ws.Range('E38').Value = '  +2.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05155'
$ws.Range('E39').Value = '  +1.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5167'
$ws.Range('E40').Value = '  +1.93%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.889'
$ws.Range('E41').Value = '  +1.93%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1510'
$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.108'
$ws.Range('E43').Value = '  +0.76%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.51'
$ws.Range('E44').Value = '  +5.78%  '

$ws.Range('E45').Value = '  +0.49%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4667'
$ws.Range('E46').Value = '  -0.66%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '100.26'
$ws.Range('E47').Value = '  +0.20%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.596'
$ws.Range('E48').Value = '  +1.36%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.66'
$ws.Range('E49').Value = '  +1.71%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06026'
$ws.Range('E50').Value = '  +0.75%  '

$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8769'
$ws.Range('E51').Value = '  +4.08%  '
